$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New quiz column: Q05 (shared string header in K1) ---
$ws.Range("K1").Value = "Q05"

# --- Per-student raw score edits (columns H/I) ---
$ws.Range("I4").Value = 10
$ws.Range("H5").Value = 10
$ws.Range("I5").Value = 10
$ws.Range("H13").Value = 8.5
$ws.Range("I13").Value = 10

# --- Quiz-average column J: denominator changed from 18 to 17 ---
$ws.Range("J2").Formula = "=(10/17)*10"
$ws.Range("J3").Formula = "=(9/17)*10"
$ws.Range("J4").Formula = "=(11/17)*10"
$ws.Range("J5").Formula = "=(10/17)*10"
$ws.Range("J6").Formula = "=(11/17)*10"
$ws.Range("J7").Formula = "=(13/17)*10"
$ws.Range("J11").Formula = "=(11/17)*10"
$ws.Range("J13").Formula = "=(11/17)*10"
$ws.Range("J14").Formula = "=(5/17)*10"

# --- New quiz (Q05) column K: new formulas / values ---
$ws.Range("K2").Formula = "=(24/35)*10"
$ws.Range("K3").Formula = "=(23/35)*10"
$ws.Range("K4").Formula = "=(29/35)*10"
$ws.Range("K5").Formula = "=(23/35)*10"
$ws.Range("K6").Formula = "=(32/35)*10"
$ws.Range("K7").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("K14").Formula = "=(14/35)*10"
$ws.Range("K15").Value = 0

# --- HomeWork 11 (Extra point column, rows 22-35) ---
$ws.Range("C22").Value = 8
$ws.Range("C23").Value = 7
$ws.Range("C24").Value = 9
$ws.Range("C25").Value = 8
$ws.Range("C26").Value = 10
$ws.Range("C27").Value = 9.5
$ws.Range("C29").Value = 3
$ws.Range("C30").Value = 4
$ws.Range("C31").Value = 7
$ws.Range("C32").Value = 8
$ws.Range("C33").Value = 9
$ws.Range("C34").Value = 6
$ws.Range("C35").Value = 7

# --- Restore selection to I5 (matches author's last cursor position) ---
$ws.Range("I5").Select()
